$d = $word.ActiveDocument
$d.Content.Find.Execute("852÷2=426, 0", $true, $false, $false, $false, $false, $true, 1, $false, "752÷3=250, 2", 2) | Out-Null
$d.Content.Find.Execute("525÷5=105, 0", $true, $false, $false, $false, $false, $true, 1, $false, "346÷4=86, 2", 2) | Out-Null
$d.Content.Find.Execute("890÷6=148, 2", $true, $false, $false, $false, $false, $true, 1, $false, "722÷8=90, 2", 2) | Out-Null
$d.Content.Find.Execute("486÷5=97, 1", $true, $false, $false, $false, $false, $true, 1, $false, "587÷8=73, 3", 2) | Out-Null
$d.Content.Find.Execute("647÷8=80, 7", $true, $false, $false, $false, $false, $true, 1, $false, "832÷9=92, 4", 2) | Out-Null
$d.Content.Find.Execute("898÷9=99, 7", $true, $false, $false, $false, $false, $true, 1, $false, "159÷2=79, 1", 2) | Out-Null
$d.Content.Find.Execute("245÷5=49, 0", $true, $false, $false, $false, $false, $true, 1, $false, "249÷9=27, 6", 2) | Out-Null
$d.Content.Find.Execute("330÷4=82, 2", $true, $false, $false, $false, $false, $true, 1, $false, "865÷9=96, 1", 2) | Out-Null
$d.Content.Find.Execute("706÷5=141, 1", $true, $false, $false, $false, $false, $true, 1, $false, "629÷8=78, 5", 2) | Out-Null
$d.Content.Find.Execute("586÷5=117, 1", $true, $false, $false, $false, $false, $true, 1, $false, "880÷4=220, 0", 2) | Out-Null
$d.Content.Find.Execute("265÷9=29, 4", $true, $false, $false, $false, $false, $true, 1, $false, "418÷3=139, 1", 2) | Out-Null
$d.Content.Find.Execute("439÷4=109, 3", $true, $false, $false, $false, $false, $true, 1, $false, "978÷5=195, 3", 2) | Out-Null
$d.Content.Find.Execute("346÷8=43, 2", $true, $false, $false, $false, $false, $true, 1, $false, "375÷7=53, 4", 2) | Out-Null
$d.Content.Find.Execute("605÷6=100, 5", $true, $false, $false, $false, $false, $true, 1, $false, "633÷7=90, 3", 2) | Out-Null
$d.Content.Find.Execute("424÷2=212, 0", $true, $false, $false, $false, $false, $true, 1, $false, "706÷4=176, 2", 2) | Out-Null
$d.Content.Find.Execute("261÷7=37, 2", $true, $false, $false, $false, $false, $true, 1, $false, "710÷8=88, 6", 2) | Out-Null
$d.Content.Find.Execute("829÷5=165, 4", $true, $false, $false, $false, $false, $true, 1, $false, "494÷5=98, 4", 2) | Out-Null
$d.Content.Find.Execute("200÷6=33, 2", $true, $false, $false, $false, $false, $true, 1, $false, "509÷5=101, 4", 2) | Out-Null
$d.Content.Find.Execute("948÷7=135, 3", $true, $false, $false, $false, $false, $true, 1, $false, "385÷5=77, 0", 2) | Out-Null
$d.Content.Find.Execute("809÷2=404, 1", $true, $false, $false, $false, $false, $true, 1, $false, "520÷7=74, 2", 2) | Out-Null
$d.Content.Find.Execute("206÷4=51, 2", $true, $false, $false, $false, $false, $true, 1, $false, "628÷4=157, 0", 2) | Out-Null
$d.Content.Find.Execute("812÷8=101, 4", $true, $false, $false, $false, $false, $true, 1, $false, "494÷6=82, 2", 2) | Out-Null
$d.Content.Find.Execute("800÷4=200, 0", $true, $false, $false, $false, $false, $true, 1, $false, "844÷9=93, 7", 2) | Out-Null
$d.Content.Find.Execute("433÷6=72, 1", $true, $false, $false, $false, $false, $true, 1, $false, "429÷5=85, 4", 2) | Out-Null
$d.Content.Find.Execute("826÷5=165, 1", $true, $false, $false, $false, $false, $true, 1, $false, "250÷9=27, 7", 2) | Out-Null
